$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly refined HGF_305-470 data (row 33, columns C:L) ---
$ws.Range("C33").Value = 50.06
$ws.Range("D33").Value = 4
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 2.79
$ws.Range("G33").Value = 32
$ws.Range("H33").Value = 4.32
$ws.Range("I33").Value = 96
$ws.Range("J33").Value = 0.0067
$ws.Range("K33").Value = 1.85
$ws.Range("L33").Value = 83

# Recalculate so the row-74 column AVERAGE formulas pick up the new row 33 data
$excel.CalculateFull()

# --- Update the view: scroll position / active selection moved to N61 ---
$ws.Range("N61").Select()
